$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: update K14 (value changed due to addition of a new diagonal/rolling point)
$ws.Range("K14").Value = 1.388747888886706

# Row 15: update J15 and add new K15
$ws.Range("J15").Value = 0.444773652920949
$ws.Range("K15").Value = 0.2348700177716323

# Row 16: update I16 and add new J16
$ws.Range("I16").Value = 0.4487415504340581
$ws.Range("J16").Value = 0.2388379152847414

# Row 17: update H17 and add new I17
$ws.Range("H17").Value = 0.5843816406042994
$ws.Range("I17").Value = 0.3744780054549828

# Row 18: update G18 and add new H18
$ws.Range("G18").Value = 0.3435754587486348
$ws.Range("H18").Value = 0.1336718235993181

# Row 19: update F19 and add new G19
$ws.Range("F19").Value = 0.2982442434965384
$ws.Range("G19").Value = 0.08834060834722172

# Row 20: update E20 and add new F20
$ws.Range("E20").Value = 0.2313828215604846
$ws.Range("F20").Value = 0.02147918641116785

# Row 21: update D21 and add new E21
$ws.Range("D21").Value = 0.201796619203768
$ws.Range("E21").Value = -0.00810701594554874

# Row 22: update C22 and add new D22
$ws.Range("C22").Value = 0.1836459624741271
$ws.Range("D22").Value = -0.02625767267518964

# Row 23: update B23 and add new C23
$ws.Range("B23").Value = 0.1656141382254278
$ws.Range("C23").Value = -0.04428949692388896

# Row 24: add new B24
$ws.Range("B24").Value = -0.09587373626955231
